# Add two new columns, I ("I0") and J ("IF"), to the active sheet,
# matching the header style already used by the other header cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy the style from the existing header cell H1
# so the new headers look consistent with the rest of the row (reuses
# the same style index rather than minting a new one).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-18 for columns I and J.
$data = @(
    @(1, 1),
    @(1, 2),
    @(8, 8),
    @(1, 1),
    @(8, 8),
    @(8, 8),
    @(10, 10),
    @(6, 7),
    @(7, 7),
    @(8, 9),
    @(7, 8),
    @(1, 1),
    @(8, 8),
    @(5, 6),
    @(9, 9),
    @(7, 8),
    @(6, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
